$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed snapshot timestamp (07:50 -> 08:20)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 08:20"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 142735
$ws.Range("C4").Value = 275
$ws.Range("D4").Value = 4559
$ws.Range("E4").Value = 135688
$ws.Range("F4").Value = 2970
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 2488

# Row 24: Suecia
$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 3700
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 3574
$ws.Range("F24").Value = 281
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 110

# Row 35: Pakistan
$ws.Range("A35").Value = "Pakistan"
$ws.Range("B35").Value = 1625
$ws.Range("C35").Value = 28
$ws.Range("D35").Value = 29
$ws.Range("E35").Value = 1578
$ws.Range("F35").Value = 11
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 18

# Row 72: Bulgaria
$ws.Range("A72").Value = "Bulgaria"
$ws.Range("B72").Value = 354
$ws.Range("C72").Value = 8
$ws.Range("D72").Value = 15
$ws.Range("E72").Value = 331
$ws.Range("F72").Value = 13
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 8

# Row 73: Letonia
$ws.Range("A73").Value = "Letonia"
$ws.Range("B73").Value = 347
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 346
$ws.Range("F73").Value = 3
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 8

# Row 100: Honduras
$ws.Range("A100").Value = "Honduras"
$ws.Range("B100").Value = 139
$ws.Range("C100").Value = 29
$ws.Range("D100").Value = 3
$ws.Range("E100").Value = 133
$ws.Range("F100").Value = 4
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 3

# Row 101: Cuba
$ws.Range("A101").Value = "Cuba"
$ws.Range("B101").Value = 139
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 4
$ws.Range("E101").Value = 132
$ws.Range("F101").Value = 2
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 3

# Row 102: Camerun
$ws.Range("A102").Value = "Camerun"
$ws.Range("B102").Value = 139
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 5
$ws.Range("E102").Value = 128
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 6

# Row 103: Brunei
$ws.Range("A103").Value = "Brunei"
$ws.Range("B103").Value = 126
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 34
$ws.Range("E103").Value = 91
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 1

# Row 104: Afganistan
$ws.Range("A104").Value = "Afganistan"
$ws.Range("B104").Value = 120
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 2
$ws.Range("E104").Value = 114
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 4

# Row 105: Sri Lanka
$ws.Range("A105").Value = "Sri Lanka"
$ws.Range("B105").Value = 120
$ws.Range("C105").Value = 3
$ws.Range("D105").Value = 11
$ws.Range("E105").Value = 108
$ws.Range("F105").Value = 5
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 1

# Row 106: Venezuela
$ws.Range("A106").Value = "Venezuela"
$ws.Range("B106").Value = 119
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 39
$ws.Range("E106").Value = 77
$ws.Range("F106").Value = 6
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 3

# Row 107: Nigeria
$ws.Range("A107").Value = "Nigeria"
$ws.Range("B107").Value = 111
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 3
$ws.Range("E107").Value = 107
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 1

# Row 108: Mauricio
$ws.Range("A108").Value = "Mauricio"
$ws.Range("B108").Value = 110
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 107
$ws.Range("F108").Value = 1
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 3

# Row 135: Uganda
$ws.Range("A135").Value = "Uganda"
$ws.Range("B135").Value = 33
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 33
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

# Row 136: Barbados
$ws.Range("A136").Value = "Barbados"
$ws.Range("B136").Value = 33
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 0
$ws.Range("E136").Value = 33
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 158: Islas Caimanes
$ws.Range("A158").Value = "Islas Caimanes"
$ws.Range("B158").Value = 12
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 11
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 1

# Row 159: Mongolia
$ws.Range("A159").Value = "Mongolia"
$ws.Range("B159").Value = 12
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 2
$ws.Range("E159").Value = 10
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0
